# Update: Threat Alert Report - 2026-01-11 01:17
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 date. Assigning the plain string would let Excel infer a
# date and convert the cell to a date serial, so enter it with a leading
# apostrophe (forces text, same as typing it in the Excel UI) and then
# restore the regular cell formatting from a neighboring data cell.
$ws.Range("A2").Value = "'16-JAN-26"
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A2").PasteSpecial(-4122) | Out-Null

# Update row 4 date and fare figures
$ws.Range("A4").Value = "'06-FEB-26"
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("E4").Value = 519
$ws.Range("F4").Value = -22

# Update row 5 fare figures
$ws.Range("E5").Value = 519
$ws.Range("F5").Value = -22

# Remove rows 6 through 10 (no longer part of the report)
$ws.Range("A6:K10").EntireRow.Delete()

# Narrow the IMPACT column (J) now that long rows are gone
$ws.Columns.Item(10).ColumnWidth = 11.17
